# Apply latest crypto price/volume snapshot (GitHub Actions update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.350.98"
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("D3").Value = "2.606.56"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'542.31"
$ws.Range("E5").Value = "  +4.33%  "
$ws.Range("D6").Value = "'141.79"
$ws.Range("E6").Value = "  +1.86%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("E11").Value = "  +2.03%  "
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "3.060.24"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").Value = "59.291.87"
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").Value = "'20.65"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000134"
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.598.73"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "'341.61"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").Value = "'4.36"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").Value = "'10.16"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").Value = "'6.41"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "'67.61"
$ws.Range("E23").Value = "  +2.04%  "
$ws.Range("D24").Value = "'0.410"
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("D25").Value = "'0.165"
$ws.Range("E25").Value = "  -1.68%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'7.26"
$ws.Range("E27").Value = "  +3.35%  "
$ws.Range("D28").Value = "0.0₃0748"
$ws.Range("E28").Value = "  +4.43%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "'1.69"
$ws.Range("E30").Value = "  +7.13%  "
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("D32").Value = "'18.77"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").Value = "'149.63"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("D34").Value = "'3.99"
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").Value = "'37.24"
$ws.Range("E36").Value = "  +1.94%  "
$ws.Range("D37").Value = "'1.47"
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("D38").Value = "'0.838"
$ws.Range("E38").Value = "  +1.57%  "
$ws.Range("E39").Value = "  +1.95%  "
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").Value = "'276.16"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").Value = "'0.600"
$ws.Range("E43").Value = "  +1.91%  "
$ws.Range("D45").Value = "'0.0956"
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("D46").Value = "'0.0525"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").Value = "1.954.53"
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'18.59"
$ws.Range("E48").Value = "  +4.16%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0224"
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("D51").Value = "'111.03"
$ws.Range("E51").Value = "  -0.48%  "
